$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(866).Insert()
$ws.Range("A866").NumberFormat = "@"
$ws.Range("A866").Value = "2026/02/27"
$ws.Range("A866").ClearFormats()
$ws.Range("B866").Value = "金"
$ws.Range("C866").Value = 1
$ws.Range("D866").Value = 32
